$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 115
$ws.Range("H115").Value = 1894.5
$ws.Range("I115").Value = 585
$ws.Range("J115").Value = 3204
$ws.Range("K115").Value = 1755
$ws.Range("L115").Value = 9612
$ws.Range("M115").Value = -188
$ws.Range("N115").Value = -12746

# Row 137
$ws.Range("H137").Value = 2704422.8
$ws.Range("I137").Value = 3334681.5
$ws.Range("J137").Value = 3314
$ws.Range("K137").Value = 10004044.5
$ws.Range("L137").Value = 9942
$ws.Range("M137").Value = -10001494.5

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 1600.4286
$ws.Range("I2").Value = 1283.8334
$ws.Range("J2").Value = 3500
$ws.Range("K2").Value = 1283.8334
$ws.Range("L2").Value = 3500
$ws.Range("M2").Value = -1170.8334
$ws.Range("N2").Value = -3726

# Row 61
$ws.Range("H61").Value = 58942388
$ws.Range("I61").Value = 77000910
$ws.Range("J61").Value = 252175
$ws.Range("K61").Value = 77000910
$ws.Range("L61").Value = 252175
$ws.Range("M61").Value = -77000698
$ws.Range("N61").Value = -252599

# Row 74
$ws.Range("H74").Value = 11455587
$ws.Range("I74").Value = 15688472
$ws.Range("J74").Value = 167892.33
$ws.Range("K74").Value = 15688472
$ws.Range("L74").Value = 167892.33
$ws.Range("M74").Value = -15687598
$ws.Range("N74").Value = -169640.33

# Row 77
$ws.Range("H77").Value = 11455587
$ws.Range("I77").Value = 15688472
$ws.Range("J77").Value = 167892.33
$ws.Range("K77").Value = 78442360
$ws.Range("L77").Value = 839461.6499999999
$ws.Range("M77").Value = -78437992
$ws.Range("N77").Value = -848197.6499999999

# Row 97
$ws.Range("H97").Value = 2402.5
$ws.Range("I97").Value = 2870
$ws.Range("J97").Value = 1000
$ws.Range("K97").Value = 2870
$ws.Range("L97").Value = 1000
$ws.Range("M97").Value = -2374

# Row 102
$ws.Range("H102").Value = 17858194
$ws.Range("I102").Value = 17858194
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 17858194
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -17856572

# Row 110
$ws.Range("H110").Value = 1123921.4
$ws.Range("I110").Value = 1444756
$ws.Range("J110").Value = 1000
$ws.Range("K110").Value = 1444756
$ws.Range("L110").Value = 1000
$ws.Range("M110").Value = -1442711
$ws.Range("N110").Value = -5090

# Row 116
$ws.Range("H116").Value = 1600.4286
$ws.Range("I116").Value = 1283.8334
$ws.Range("J116").Value = 3500
$ws.Range("K116").Value = 1283.8334
$ws.Range("L116").Value = 3500
$ws.Range("M116").Value = 1010.1666
$ws.Range("N116").Value = -8088

# Row 122
$ws.Range("H122").Value = 1318.7142
$ws.Range("I122").Value = 1149.5333
$ws.Range("J122").Value = 1741.6666
$ws.Range("K122").Value = 3448.5999
$ws.Range("L122").Value = 5224.9998
$ws.Range("M122").Value = -998.5999000000002
$ws.Range("N122").Value = -10124.9998

# Row 132
$ws.Range("H132").Value = 36562.9
$ws.Range("I132").Value = 24469.256
$ws.Range("J132").Value = 67152.7
$ws.Range("K132").Value = 73407.76800000001
$ws.Range("L132").Value = 201458.1
$ws.Range("M132").Value = -70877.76800000001
$ws.Range("N132").Value = -206518.1

# Row 136
$ws.Range("H136").Value = 58942388
$ws.Range("I136").Value = 77000910
$ws.Range("J136").Value = 252175
$ws.Range("K136").Value = 231002730
$ws.Range("L136").Value = 756525
$ws.Range("M136").Value = -231000180
$ws.Range("N136").Value = -761625

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 1600.4286
$ws.Range("I3").Value = 1283.8334
$ws.Range("J3").Value = 3500
$ws.Range("K3").Value = 1283.8334
$ws.Range("L3").Value = 3500
$ws.Range("M3").Value = -1169.8334
$ws.Range("N3").Value = -3728

# Row 80
$ws.Range("H80").Value = 483.66666
$ws.Range("I80").Value = 118
$ws.Range("J80").Value = 539.9231
$ws.Range("K80").Value = 118
$ws.Range("L80").Value = 539.9231
$ws.Range("M80").Value = 880
$ws.Range("N80").Value = -2535.9231

# Row 83
$ws.Range("H83").Value = 483.66666
$ws.Range("I83").Value = 118
$ws.Range("J83").Value = 539.9231
$ws.Range("K83").Value = 590
$ws.Range("L83").Value = 2699.6155
$ws.Range("M83").Value = 4402
$ws.Range("N83").Value = -12683.6155

# Row 94
$ws.Range("H94").Value = 837.4545000000001
$ws.Range("I94").Value = 455.84616
$ws.Range("J94").Value = 1388.6666
$ws.Range("K94").Value = 455.84616
$ws.Range("L94").Value = 1388.6666
$ws.Range("M94").Value = -4.846159999999998
$ws.Range("N94").Value = -2290.6666

# Row 99
$ws.Range("H99").Value = 1657
$ws.Range("I99").Value = 810
$ws.Range("J99").Value = 1826.4
$ws.Range("K99").Value = 810
$ws.Range("L99").Value = 1826.4
$ws.Range("M99").Value = 688
$ws.Range("N99").Value = -4822.4

# Row 107
$ws.Range("H107").Value = 2688.4614
$ws.Range("I107").Value = 2221.5715
$ws.Range("J107").Value = 3233.1667
$ws.Range("K107").Value = 2221.5715
$ws.Range("L107").Value = 3233.1667
$ws.Range("M107").Value = -301.5715
$ws.Range("N107").Value = -7073.1667

# Row 134
$ws.Range("H134").Value = 1522.6487
$ws.Range("I134").Value = 1574.963
$ws.Range("J134").Value = 1381.4
$ws.Range("K134").Value = 4724.889
$ws.Range("L134").Value = 4144.200000000001
$ws.Range("M134").Value = -2189.889

$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 41668944
$ws.Range("I16").Value = 1550
$ws.Range("J16").Value = 83336340
$ws.Range("K16").Value = 1550
$ws.Range("L16").Value = 83336340
$ws.Range("M16").Value = -1263
$ws.Range("N16").Value = -83336914

# Row 22
$ws.Range("H22").Value = 557.8889
$ws.Range("I22").Value = 264.2
$ws.Range("J22").Value = 925
$ws.Range("K22").Value = 264.2
$ws.Range("L22").Value = 925
$ws.Range("M22").Value = 85.80000000000001
$ws.Range("N22").Value = -1625

# Row 31
$ws.Range("H31").Value = 2502.3489
$ws.Range("I31").Value = 1230.9395
$ws.Range("J31").Value = 6698
$ws.Range("K31").Value = 1230.9395
$ws.Range("L31").Value = 6698
$ws.Range("M31").Value = -935.9395

# Row 34
$ws.Range("H34").Value = 2502.3489
$ws.Range("I34").Value = 1230.9395
$ws.Range("J34").Value = 6698
$ws.Range("K34").Value = 1230.9395
$ws.Range("L34").Value = 6698
$ws.Range("M34").Value = -1028.9395

# Row 58
$ws.Range("H58").Value = 31254014
$ws.Range("I58").Value = 45458290
$ws.Range("J58").Value = 4609.5
$ws.Range("K58").Value = 45458290
$ws.Range("L58").Value = 4609.5
$ws.Range("M58").Value = -45458087

# Row 99
$ws.Range("H99").Value = 2627.4
$ws.Range("I99").Value = 2517.625
$ws.Range("J99").Value = 2822.5557
$ws.Range("K99").Value = 2517.625
$ws.Range("L99").Value = 2822.5557
$ws.Range("M99").Value = -1019.625
$ws.Range("N99").Value = -5818.5557

# Row 105
$ws.Range("H105").Value = 1767.8889
$ws.Range("I105").Value = 600
$ws.Range("J105").Value = 2702.2
$ws.Range("K105").Value = 600
$ws.Range("L105").Value = 2702.2
$ws.Range("M105").Value = 1147
$ws.Range("N105").Value = -6196.2

# Row 113
$ws.Range("H113").Value = 41668944
$ws.Range("I113").Value = 1550
$ws.Range("J113").Value = 83336340
$ws.Range("K113").Value = 1550
$ws.Range("L113").Value = 83336340
$ws.Range("M113").Value = 620
$ws.Range("N113").Value = -83340680

# Row 122
$ws.Range("H122").Value = 2148.7932
$ws.Range("I122").Value = 1639.4445
$ws.Range("J122").Value = 2982.2727
$ws.Range("K122").Value = 4918.333500000001
$ws.Range("L122").Value = 8946.8181
$ws.Range("M122").Value = -2468.333500000001
$ws.Range("N122").Value = -13846.8181

# Row 126
$ws.Range("H126").Value = 2627.4
$ws.Range("I126").Value = 2517.625
$ws.Range("J126").Value = 2822.5557
$ws.Range("K126").Value = 7552.875
$ws.Range("L126").Value = 8467.667099999999
$ws.Range("M126").Value = -5082.875
$ws.Range("N126").Value = -13407.6671

# Row 132
$ws.Range("H132").Value = 26319.17
$ws.Range("I132").Value = 1767.35
$ws.Range("J132").Value = 49701.855
$ws.Range("K132").Value = 5302.049999999999
$ws.Range("L132").Value = 149105.565
$ws.Range("M132").Value = -2772.049999999999
$ws.Range("N132").Value = -154165.565

# Row 134
$ws.Range("H134").Value = 35542.688
$ws.Range("I134").Value = 1903.68
$ws.Range("J134").Value = 155682
$ws.Range("K134").Value = 5711.04
$ws.Range("L134").Value = 467046
$ws.Range("M134").Value = -3176.04

# Row 136
$ws.Range("H136").Value = 31254014
$ws.Range("I136").Value = 45458290
$ws.Range("J136").Value = 4609.5
$ws.Range("K136").Value = 136374870
$ws.Range("L136").Value = 13828.5
$ws.Range("M136").Value = -136372320

$ws = $wb.Worksheets.Item("CUL")
# Row 131
$ws.Range("H131").Value = 818.4666999999999
$ws.Range("I131").Value = 635.7143
$ws.Range("J131").Value = 978.375
$ws.Range("K131").Value = 1907.1429
$ws.Range("L131").Value = 2935.125
$ws.Range("M131").Value = 3132.8571
$ws.Range("N131").Value = -13015.125

$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 4483
$ws.Range("I70").Value = 4142.9473
$ws.Range("J70").Value = 5200.8887
$ws.Range("K70").Value = 4142.9473
$ws.Range("L70").Value = 5200.8887
$ws.Range("M70").Value = -3872.9473
$ws.Range("N70").Value = -5740.8887

# Row 73
$ws.Range("H73").Value = 4483
$ws.Range("I73").Value = 4142.9473
$ws.Range("J73").Value = 5200.8887
$ws.Range("K73").Value = 4142.9473
$ws.Range("L73").Value = 5200.8887
$ws.Range("M73").Value = -3206.9473
$ws.Range("N73").Value = -7072.8887

# Row 102
$ws.Range("H102").Value = 1383
$ws.Range("I102").Value = 1356.7778
$ws.Range("J102").Value = 1501
$ws.Range("K102").Value = 1356.7778
$ws.Range("L102").Value = 1501
$ws.Range("M102").Value = 265.2221999999999

# Row 113
$ws.Range("H113").Value = 1745.7727
$ws.Range("I113").Value = 1622.6111
$ws.Range("J113").Value = 2300
$ws.Range("K113").Value = 1622.6111
$ws.Range("L113").Value = 2300
$ws.Range("M113").Value = 547.3888999999999
$ws.Range("N113").Value = -6640

# Row 126
$ws.Range("H126").Value = 1933.3334
$ws.Range("I126").Value = 1866.6666
$ws.Range("J126").Value = 2000
$ws.Range("K126").Value = 5599.9998
$ws.Range("L126").Value = 6000
$ws.Range("M126").Value = -3129.9998

# Row 132
$ws.Range("H132").Value = 67902.836
$ws.Range("I132").Value = 44249
$ws.Range("J132").Value = 145622.58
$ws.Range("K132").Value = 132747
$ws.Range("L132").Value = 436867.74
$ws.Range("M132").Value = -130217

$ws = $wb.Worksheets.Item("LTW")
# Row 132
$ws.Range("H132").Value = 37829.55
$ws.Range("I132").Value = 2474.4707
$ws.Range("J132").Value = 87915.914
$ws.Range("K132").Value = 7423.4121
$ws.Range("L132").Value = 263747.742
$ws.Range("M132").Value = -4893.4121

# Row 136
$ws.Range("H136").Value = 183443.1
$ws.Range("I136").Value = 112876
$ws.Range("J136").Value = 500995
$ws.Range("K136").Value = 338628
$ws.Range("L136").Value = 1502985
$ws.Range("M136").Value = -336078

$ws = $wb.Worksheets.Item("WVR")
# Row 107
$ws.Range("H107").Value = 217.33333
$ws.Range("I107").Value = 220.8
$ws.Range("J107").Value = 200
$ws.Range("K107").Value = 662.4000000000001
$ws.Range("L107").Value = 600
$ws.Range("M107").Value = 1257.6
$ws.Range("N107").Value = -4440

# Row 132
$ws.Range("H132").Value = 57126.36
$ws.Range("I132").Value = 37922.816
$ws.Range("J132").Value = 114737
$ws.Range("K132").Value = 113768.448
$ws.Range("L132").Value = 344211
$ws.Range("M132").Value = -111238.448

# Row 136
$ws.Range("H136").Value = 65481.484
$ws.Range("I136").Value = 42692.75
$ws.Range("J136").Value = 143614.28
$ws.Range("K136").Value = 128078.25
$ws.Range("L136").Value = 430842.84
$ws.Range("M136").Value = -125528.25
$ws.Range("N136").Value = -435942.84
